# Mise à jour des activité et du cours Excel
#
# - Rename the two worksheets ("Activité 10" -> "Activité 14",
#   "Activité 11" -> "Activité 15").
# - Move the active tab from the first sheet to the second one
#   (workbook activeTab, and each sheet's tabSelected flag).
# - Update the header/footer font style from "Times New Roman,Regular"
#   to "Times New Roman,Normal" on every sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename sheets.
$ws1.Name = "Activité 14"
$ws2.Name = "Activité 15"

# Make the second sheet the active / selected tab.
$ws2.Activate()

# Update the header & footer font name on both sheets.
foreach ($ws in @($ws1, $ws2)) {
    $ws.PageSetup.CenterHeader = '&"Times New Roman,Normal"&12&A'
    $ws.PageSetup.CenterFooter = '&"Times New Roman,Normal"&12Page &P'
}
